# "Mas tratamiento de datos - Incluido el one-hot"
# Add one-hot encoded consensus columns (concensoV1.1 / concensoV1.2) to the
# metrics sheet, filling the new cells with the placeholder value "~".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill the new one-hot placeholder values first so the shared-string table
# picks up "~" before the new header labels.
$ws.Range("C2:D5").Value = "~"

# Rename the existing consensus header and add the new one-hot headers.
$ws.Range("C1").Value = "concensoV1.1"

# New column D: copy C1's formatting (bold, centered header style) onto D1,
# then set its header text.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "concensoV1.2"

# Give the new column a slightly wider width to fit the longer header.
$ws.Columns.Item(4).ColumnWidth = 12.3

# Match the saved selection/active cell state.
$ws.Range("D2").Select()
